# sprint3 final DOCs update
# Update the Sprint 3 burn-down tracker (Sheet1) with the final totals:
#  - Backlog (C2) grew from 37 to 45 points
#  - Completed-work entries recorded for the last few days of the sprint
#    (D14:D17), which were previously blank
# All dependent formulas (E - BurnDn, F - Ideal, the C26/D26 sums, and the
# chart series that read from Sheet1!E2:E25 / Sheet1!F2:F25) recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Total backlog size for the sprint increased from 37 to 45
$ws.Range("C2").Value = 45

# Newly-logged "Completed" work for days 13-16 of the sprint
$ws.Range("D14").Value = 5
$ws.Range("D15").Value = 4
$ws.Range("D16").Value = 18
$ws.Range("D17").Value = 11

# Force a full recalculation so BurnDn/Ideal columns, the sums, and the
# chart caches are all refreshed before the workbook is saved
$excel.CalculateFullRebuild()

# Update the sheet's active cell/selection
[void]$ws.Range("W4").Select()
